$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.764.52"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.488.49"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'588.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'169.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("D8").Value = "3.477.41"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").Value = "'0.573"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.38%  "
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'0.0000276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "4.053.60"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'616.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.78%  "
$ws.Range("E17").Value = "  -4.87%  "
$ws.Range("D18").Value = "3.476.63"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "68.783.05"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").Value = "'17.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "'11.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("D24").Value = "'15.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("D25").Value = "'95.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").Value = "'9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("D30").Value = "'32.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "'8.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.93%  "
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "'1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  -6.14%  "
$ws.Range("D35").Value = "'577.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'3.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.95%  "
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "'57.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").Value = "'0.0436"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "3.413.84"
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("D45").Value = "'32.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "0.0₃0695"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "'2.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("E50").Value = "  +13.30%  "
$ws.Range("D51").Value = "'132.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.34%  "
